$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Department/Faculty for the last record (row 4)
$ws.Range("D4").Value = "ICT"
$ws.Range("E4").Value = "Basic and Applied Sciences"

# Move the active selection to E5 (next empty row under the table)
$ws.Range("E5").Select()
